$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diff effectively swaps the values in columns D, M, N, O, P, S between
# row 2 and row 3 (all other columns are identical between the two rows).

# Row 2 (new values, formerly row 3's values)
$ws.Range("D2").Value = 44209
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 11000
$ws.Range("P2").Value = 10500
$ws.Range("S2").Value = 750

# Row 3 (new values, formerly row 2's values)
$ws.Range("D3").Value = 44217
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11500
$ws.Range("S3").Value = 821
